$wb = $excel.ActiveWorkbook

$wsArbeitspakete   = $wb.Worksheets.Item(1)   # "Arbeitspakete"
$wsArbeitspaketePM = $wb.Worksheets.Item(2)   # "Arbeitspakete PM"
$wsAufwaende       = $wb.Worksheets.Item(3)   # "Aufwände gesamt"

# ---------------------------------------------------------------------
# Data fix: the "Aufwände gesamt" sheet imports the budget name from the
# wrong column. Column E already holds the budget name (shared strings
# "Budget1"/"Budget2") - column F must show the same value for every
# data row (rows 4-11).
# ---------------------------------------------------------------------
for ($r = 4; $r -le 11; $r++) {
    $wsAufwaende.Cells.Item($r, 6).Value = $wsAufwaende.Cells.Item($r, 5).Value2
}

# ---------------------------------------------------------------------
# Print titles on "Aufwände gesamt" only repeat the first three rows
# (header row + column header row) instead of four.
# ---------------------------------------------------------------------
$wsAufwaende.PageSetup.PrintTitleRows = '$1:$3'

# ---------------------------------------------------------------------
# Restore the selections as left behind after editing column F on the
# "Aufwände gesamt" sheet: the new values occupy F4:F11 there, while the
# other two sheets keep their previous anchor cell.
# ---------------------------------------------------------------------
[void]$wsArbeitspakete.Activate()
[void]$wsArbeitspakete.Range("A1").Select()

[void]$wsArbeitspaketePM.Activate()
[void]$wsArbeitspaketePM.Range("D21").Select()

[void]$wsAufwaende.Activate()
[void]$wsAufwaende.Range("F4:F11").Select()
